$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) from row 16 down into the two brand-new rows (17, 18)
# so they get the same style index as the rest of the data rows, rather than
# creating brand-new style entries.
$ws.Cells.Item(16,1).Copy()
$ws.Range("A17:B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 18: Immutable / BooleanType[null]  (previously occupied row 16)
$ws.Cells.Item(18,1).Value = "Immutable"
$ws.Cells.Item(18,2).Value = "BooleanType[null]"

# Row 17: Copyright / (empty)  (previously occupied row 15)
$ws.Cells.Item(17,1).Value = "Copyright"
$ws.Cells.Item(17,2).ClearContents()

# Row 16: Purpose / (empty)  (previously occupied row 14)
$ws.Cells.Item(16,1).Value = "Purpose"
$ws.Cells.Item(16,2).ClearContents()

# Row 15: Description / ...  (previously occupied row 13)
$ws.Cells.Item(15,1).Value = "Description"
$ws.Cells.Item(15,2).Value = "Definiert Dokumentanforderungkategorien"

# Row 14: Jurisdiction / Germany  (previously occupied row 12)
$ws.Cells.Item(14,1).Value = "Jurisdiction"
$ws.Cells.Item(14,2).Value = "Germany"

# Row 13: Contact / No display for ContactDetail  (new contact entry)
$ws.Cells.Item(13,1).Value = "Contact"
$ws.Cells.Item(13,2).Value = "No display for ContactDetail"

# Row 12: Contact / No display for ContactDetail  (new contact entry)
$ws.Cells.Item(12,1).Value = "Contact"
$ws.Cells.Item(12,2).Value = "No display for ContactDetail"
